$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price entry was recorded. It belongs right after the existing
# entry for 2021-07-15 (row 3), so insert a blank row at position 4 and shift
# the rest of the data (previously rows 4-24) down to rows 5-25.
$ws.Rows.Item(4).Insert()

# Fill in the newly inserted row 4 with the new weekly record.
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44749
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 100112035
$ws.Range("G4").Value = "Bruselas (repollito)"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 220
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 19091
$ws.Range("N4").Value = "`$/malla 15 kilos"
$ws.Range("O4").Value = "Provincia de Quillota"
$ws.Range("P4").Value = 1273
$ws.Range("Q4").Value = 15
$ws.Range("R4").Value = "Hortaliza"
